$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text type (avoids Excel auto-numeric
# conversion for values like "1.000" or "304.50"), then restore the original
# cell style so no stray style index is introduced.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '23.800.16'
$ws.Range('E2').Value = '  +2.11%  '
$ws.Range('D3').Value = '1.653.91'
$ws.Range('E3').Value = '  +1.96%  '
Set-TextValue 'D4' '1.000'
$ws.Range('E4').Value = '  -0.13%  '
Set-TextValue 'D5' '1.000'
$ws.Range('E5').Value = '  -0.09%  '
Set-TextValue 'D6' '304.50'
$ws.Range('E6').Value = '  +0.67%  '
Set-TextValue 'D7' '0.3819'
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D8' '0.3613'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D9' '51.18'
$ws.Range('E9').Value = '  -0.42%  '
Set-TextValue 'D10' '1.250'
$ws.Range('E10').Value = '  +2.64%  '
Set-TextValue 'D11' '0.08223'
$ws.Range('E11').Value = '  +1.17%  '
Set-TextValue 'D12' '0.9996'
$ws.Range('E12').Value = '  -0.20%  '
Set-TextValue 'D13' '22.70'
$ws.Range('E13').Value = '  +2.06%  '
Set-TextValue 'D14' '6.536'
$ws.Range('E14').Value = '  +1.31%  '
Set-TextValue 'D15' '7.422'
$ws.Range('E15').Value = '  +2.33%  '
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range('D17').Value = '1.642.25'
$ws.Range('E17').Value = '  +1.63%  '
Set-TextValue 'D18' '97.79'
$ws.Range('E18').Value = '  +4.08%  '
Set-TextValue 'D19' '0.06976'
$ws.Range('E19').Value = '  +0.74%  '
Set-TextValue 'D20' '6.775'
$ws.Range('E20').Value = '  +3.92%  '
Set-TextValue 'D21' '17.78'
$ws.Range('E21').Value = '  +1.80%  '
Set-TextValue 'D22' '0.9997'
$ws.Range('E22').Value = '  -0.18%  '
Set-TextValue 'D23' '12.71'
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('D24').Value = '23.791.57'
$ws.Range('E24').Value = '  +2.10%  '
Set-TextValue 'D25' '2.542'
$ws.Range('E25').Value = '  +3.04%  '
Set-TextValue 'D26' '3.098'
$ws.Range('E26').Value = '  +0.78%  '
Set-TextValue 'D27' '21.31'
$ws.Range('E27').Value = '  +1.01%  '
Set-TextValue 'D28' '151.07'
$ws.Range('E28').Value = '  +0.36%  '
Set-TextValue 'D29' '5.240'
$ws.Range('E29').Value = '  -0.52%  '
Set-TextValue 'D30' '134.63'
$ws.Range('E30').Value = '  +1.47%  '
$ws.Range('D31').Value = '1.830.24'
$ws.Range('E31').Value = '  +1.75%  '
Set-TextValue 'D32' '6.903'
$ws.Range('E32').Value = '  +2.84%  '
Set-TextValue 'D33' '1.086'
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D34' '2.123'
$ws.Range('E34').Value = '  -1.01%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D35' '11.93'
$ws.Range('E35').Value = '  +6.24%  '
Set-TextValue 'D36' '0.02853'
$ws.Range('E36').Value = '  +4.07%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D37' '0.2522'
$ws.Range('E37').Value = '  +2.00%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D38' '6.155'
$ws.Range('E38').Value = '  +3.14%  '
Set-TextValue 'D39' '0.08835'
$ws.Range('E39').Value = '  +0.57%  '
Set-TextValue 'D40' '0.07162'
$ws.Range('E40').Value = '  +1.19%  '
Set-TextValue 'D41' '12.93'
$ws.Range('E41').Value = '  +7.89%  '
Set-TextValue 'D42' '0.7088'
$ws.Range('E42').Value = '  +1.97%  '
Set-TextValue 'D43' '1.342'
$ws.Range('E43').Value = '  +1.16%  '
Set-TextValue 'D44' '15.88'
$ws.Range('E44').Value = '  -1.04%  '
Set-TextValue 'D45' '0.6562'
$ws.Range('E45').Value = '  +1.91%  '
Set-TextValue 'D46' '2.337'
$ws.Range('E46').Value = '  +3.50%  '
Set-TextValue 'D47' '0.9995'
$ws.Range('E47').Value = '  -0.11%  '
Set-TextValue 'D48' '3.965'
$ws.Range('E48').Value = '  +0.29%  '
Set-TextValue 'D49' '0.07988'
$ws.Range('E49').Value = '  +0.41%  '
Set-TextValue 'D50' '128.90'
$ws.Range('E50').Value = '  +2.59%  '
$ws.Range('E51').Value = '  +1.42%  '
